$wb = $excel.ActiveWorkbook

# The "Date" metadata value and the "Description" metadata value both live
# as shared strings that are reused elsewhere in the workbook (the same
# "Quantité de produit" text is already used by the quantiteProduit element's
# Short/Definition columns, and the "Entrée Quantité de produit" text is also
# reused by the base element's Definition column). Doing a workbook-wide
# replace keeps every one of those references consistent with the edited
# text, exactly like the source edit did.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("2026-01-07T21:00:10+00:00", "2026-01-14T15:34:52+00:00")
    $ws.Cells.Replace("Entrée Quantité de produit", "Quantité de produit")
}
